# Apply the price-list update to "Hoja1":
#  - bump the date in A1 by one month (2024-04-24 -> 2024-05-24)
#  - apply the ~26.5% price increase to the PRECIO column entries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date cell (serial date value)
$ws.Range("A1").Value = 45436

# SOPORTE VISILLO prices
$ws.Range("D23").Value = 48.442
$ws.Range("D24").Value = 48.442
$ws.Range("D25").Value = 28.707
$ws.Range("D26").Value = 28.707

# Soporte de FLEJE prices
$ws.Range("D41").Value = 110.795
$ws.Range("D42").Value = 110.795
$ws.Range("D43").Value = 121.996
$ws.Range("D44").Value = 152.625
